# Practice on datasets for basic questions: drop the "Nan"/"None" placeholder
# text values from the weather data sheet (mirrors a pandas dropna-style cleanup)
# and leave behind truly empty cells instead of literal "Nan"/"None" strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D8").ClearContents()

# Move the active selection to D9, matching the saved cursor position.
$ws.Range("D9").Select()
